$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "O cozinheiro entrega a lista de produtos em falta"
$ws.Range("D4").Value = "O fornecedor envia o orçamento do pedido"
$ws.Range("E4").Value = "X(1)"
$ws.Range("D5").Value = "O fornecedor entrega os produtos"
$ws.Range("E5").Value = "X(2)"

$ws.Rows.Item(6).Delete()

$wb.Save()
